# Existing ICDC Biobank filter fixes
#
# The "StatQuery" (column C) used by every tab on the `startup` sheet joined
# study-level files with `OPTIONAL MATCH (sf:file)-->(s)` (file attached to
# the *study*). It should join file to *case* instead, matching the rest of
# the query: `OPTIONAL MATCH (sf:file)-->(c)`.
#
# Fix it in place for every row that carries the StatQuery (C2:C5) by doing
# a targeted text replacement on the existing (shared) string, so every
# other character -- including the non-breaking spaces sprinkled through
# the Cypher text -- is preserved byte-for-byte.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$statQueryRange = $ws.Range("C2:C5")

foreach ($cell in $statQueryRange.Cells) {
    $oldText = $cell.Text
    $newText = $oldText.Replace("OPTIONAL MATCH (sf:file)-->(s)", "OPTIONAL MATCH (sf:file)-->(c)")
    if ($newText -ne $oldText) {
        $cell.Value = $newText
    }
}

# Refresh the view to match the author's saved state: scrolled so row 2 is
# visible at the top, with the active selection on J3.
[void]$ws.Range("J3").Select()
